$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date format from the existing date cell (B11) onto the
# new rows so the new cells reuse the same style (no new style added).
$ws.Range("B11").Copy()
$ws.Range("B12:B14").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in new timesheet rows with dates (as date serials) and hour counts
$ws.Range("B12").Value2 = 44524
$ws.Range("C12").Value = 2

$ws.Range("B13").Value2 = 44525
$ws.Range("C13").Value = 3

$ws.Range("B14").Value2 = 44527

# Update the active selection to C14
$ws.Range("C14").Select()
